$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (second data row) - automatic electricity price update
$ws.Range("A2").Value = 45978
$ws.Range("B2").Value = 77.53
$ws.Range("C2").Value = 70.05
$ws.Range("D2").Value = 66.23999999999999
$ws.Range("E2").Value = 65.36
$ws.Range("F2").Value = 65.92
$ws.Range("G2").Value = 74.34
$ws.Range("H2").Value = 91.87
$ws.Range("I2").Value = 125.68
$ws.Range("J2").Value = 136.42
$ws.Range("K2").Value = 82.42
$ws.Range("L2").Value = 66.84999999999999
$ws.Range("M2").Value = 57.19
$ws.Range("N2").Value = 50.95
$ws.Range("O2").Value = 39.1
$ws.Range("P2").Value = 37.62
$ws.Range("Q2").Value = 46.51
$ws.Range("R2").Value = 69.36
$ws.Range("S2").Value = 92.48999999999999
$ws.Range("T2").Value = 112.82
$ws.Range("U2").Value = 110.74
$ws.Range("V2").Value = 103.73
$ws.Range("W2").Value = 90.19
$ws.Range("X2").Value = 87.47
$ws.Range("Y2").Value = 80.16
$ws.Range("Z2").Value = 79.20999999999999
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 96.34999999999999
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 111.78
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 109.42
$ws.Range("AG2").Value = "0h-16h"
